$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.46"
$ws.Range("D2").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05847"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.374"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.340"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8081"
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9690"
$ws.Range("D9").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07483"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03244"
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03037"
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.139"
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09401"
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001592"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04805"
$ws.Range("D17").Style = "Normal"

$ws.Range("E18").Value = "17OneONEWorstin24h"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006120"
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.004115"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0009939"
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001500"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.704"
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.229"
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3209"
$ws.Range("D25").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03873"
$ws.Range("D40").Style = "Normal"

$ws.Range("B41").Value = "KickToken"

$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006628"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"

$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1078"
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"

$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002999"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = "42CEJICEJI"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.006390"
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005611"
$ws.Range("D45").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4197"
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1465"
$ws.Range("D48").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01010"
$ws.Range("D50").Style = "Normal"
